# June 21, 2022 09:44AM
# Add three new config rows (Delay_Medium, NumberOfRetries, RetryInterval)
# to the "Assets" sheet, each tied to the "Generic Asset" asset.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")
$ws.Activate()

$newEntries = @("Delay_Medium", "NumberOfRetries", "RetryInterval")

$row = 7
foreach ($name in $newEntries) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 3).Value = "Generic Asset"
    $row++
}

$ws.Range("D9").Select()
